$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.069404619044562
$ws.Range("D2").Value = 1.074188088067269
$ws.Range("E2").Value = 1.063392055330241
$ws.Range("F2").Value = 1.082450564687021
$ws.Range("I2").Value = 1.02359499962809
$ws.Range("J2").Value = 1.074339515050223
$ws.Range("K2").Value = 1.076877962343129
$ws.Range("L2").Value = 1.066110819384368
$ws.Range("M2").Value = 1.085118749329558
$ws.Range("N2").Value = 1.075865200390059

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.073499821534699
$ws.Range("D3").Value = 1.078068119455432
$ws.Range("E3").Value = 1.067014910477174
$ws.Range("F3").Value = 1.086553944373796
$ws.Range("I3").Value = 1.023504579208684
$ws.Range("J3").Value = 1.078077137737485
$ws.Range("K3").Value = 1.08056757023439
$ws.Range("L3").Value = 1.069541647556044
$ws.Range("M3").Value = 1.089032851632673
$ws.Range("N3").Value = 1.079608130930248

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.076122029012802
$ws.Range("D4").Value = 1.080551921490242
$ws.Range("E4").Value = 1.069333299430817
$ws.Range("F4").Value = 1.089181731879248
$ws.Range("I4").Value = 1.023443335729026
$ws.Range("J4").Value = 1.080468739226115
$ws.Range("K4").Value = 1.082928205225854
$ws.Range("L4").Value = 1.07173579085114
$ws.Range("M4").Value = 1.091538249183228
$ws.Range("N4").Value = 1.082003128767313

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.077218013818548
$ws.Range("D5").Value = 1.081589905841726
$ws.Range("E5").Value = 1.070301969291543
$ws.Range("F5").Value = 1.090280129657587
$ws.Range("I5").Value = 1.023416932628352
$ws.Range("J5").Value = 1.081467945314063
$ws.Range("K5").Value = 1.08391441517548
$ws.Range("L5").Value = 1.072652223395528
$ws.Range("M5").Value = 1.092585207443028
$ws.Range("N5").Value = 1.083003753842517

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.07740166629157
$ws.Range("D6").Value = 1.081763830294934
$ws.Range("E6").Value = 1.070464268354769
$ws.Range("F6").Value = 1.090464191158561
$ws.Range("I6").Value = 1.023412460900047
$ws.Range("J6").Value = 1.081635357590053
$ws.Range("K6").Value = 1.084079646577309
$ws.Range("L6").Value = 1.072805751137322
$ws.Range("M6").Value = 1.092760632593141
$ws.Range("N6").Value = 1.083171403863141

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.076136698426811
$ws.Range("D7").Value = 1.080565815187579
$ws.Range("E7").Value = 1.069346266077438
$ws.Range("F7").Value = 1.089196433274352
$ws.Range("I7").Value = 1.02344298551034
$ws.Range("J7").Value = 1.080482114839761
$ws.Range("K7").Value = 1.082941407100334
$ws.Range("L7").Value = 1.071748059526137
$ws.Range("M7").Value = 1.091552263198383
$ws.Range("N7").Value = 1.082016523375864

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.070794481193902
$ws.Range("D8").Value = 1.075505053895057
$ws.Range("E8").Value = 1.064621892479434
$ws.Range("F8").Value = 1.083843129312902
$ws.Range("I8").Value = 1.023565007130014
$ws.Range("J8").Value = 1.075608361290324
$ws.Range("K8").Value = 1.078130559164072
$ws.Range("L8").Value = 1.067275754349051
$ws.Range("M8").Value = 1.0864473247542
$ws.Range("N8").Value = 1.077135848537359

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.061158198065165
$ws.Range("D9").Value = 1.066371616934378
$ws.Range("E9").Value = 1.056089480605245
$ws.Range("F9").Value = 1.074189616686947
$ws.Range("I9").Value = 1.023759139049229
$ws.Range("J9").Value = 1.066804264487165
$ws.Range("K9").Value = 1.069438200897041
$ws.Range("L9").Value = 1.059187887108219
$ws.Range("M9").Value = 1.077232449014184
$ws.Range("N9").Value = 1.068319248906896

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.054568734874687
$ws.Range("D10").Value = 1.060122839671844
$ws.Range("E10").Value = 1.050247859037541
$ws.Range("F10").Value = 1.06759035650928
$ws.Range("I10").Value = 1.023874603872166
$ws.Range("J10").Value = 1.060775224264197
$ws.Range("K10").Value = 1.063484407765973
$ws.Range("L10").Value = 1.053643312307951
$ws.Range("M10").Value = 1.070926757649117
$ws.Range("N10").Value = 1.062281646755271

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.051672573472101
$ws.Range("D11").Value = 1.057375686463261
$ws.Range("E11").Value = 1.047678757758961
$ws.Range("F11").Value = 1.064690386353461
$ws.Range("I11").Value = 1.023921309042974
$ws.Range("J11").Value = 1.058123332759834
$ws.Range("K11").Value = 1.060865314548612
$ws.Range("L11").Value = 1.051203099397989
$ws.Range("M11").Value = 1.068154289907685
$ws.Range("N11").Value = 1.059625989260796

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.050590041262799
$ws.Range("D12").Value = 1.05634874355818
$ws.Range("E12").Value = 1.046718233444627
$ws.Range("F12").Value = 1.063606507476276
$ws.Range("I12").Value = 1.023938164268118
$ws.Range("J12").Value = 1.057131797778757
$ws.Range("K12").Value = 1.059885997993526
$ws.Range("L12").Value = 1.050290498862824
$ws.Range("M12").Value = 1.067117839014087
$ws.Range("N12").Value = 1.058633046186314

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.050822559818213
$ws.Range("D13").Value = 1.056569326824628
$ws.Range("E13").Value = 1.04692455661153
$ws.Range("F13").Value = 1.063839311766791
$ws.Range("I13").Value = 1.023934571058996
$ws.Range("J13").Value = 1.057344784757571
$ws.Range("K13").Value = 1.060096362409735
$ws.Range("L13").Value = 1.050486539853506
$ws.Range("M13").Value = 1.067340466601683
$ws.Range("N13").Value = 1.058846335631069

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.051583231047803
$ws.Range("D14").Value = 1.057290934048184
$ws.Range("E14").Value = 1.047599489655192
$ws.Range("F14").Value = 1.064600931206835
$ws.Range("I14").Value = 1.023922712353274
$ws.Range("J14").Value = 1.058041506658047
$ws.Range("K14").Value = 1.060784497678434
$ws.Range("L14").Value = 1.051127791637372
$ws.Range("M14").Value = 1.068068753727956
$ws.Range("N14").Value = 1.059544046956559

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.052050998725618
$ws.Range("D15").Value = 1.057734665573342
$ws.Range("E15").Value = 1.048014501656274
$ws.Range("F15").Value = 1.065069292219917
$ws.Range("I15").Value = 1.023915340507297
$ws.Range("J15").Value = 1.058469908820808
$ws.Range("K15").Value = 1.061207614162277
$ws.Range("L15").Value = 1.051522058260361
$ws.Range("M15").Value = 1.068516586909232
$ws.Range("N15").Value = 1.059973057499529

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.054760012593408
$ws.Range("D16").Value = 1.060304261038092
$ws.Range("E16").Value = 1.050417502184322
$ws.Range("F16").Value = 1.067781896382619
$ws.Range("I16").Value = 1.023871434968932
$ws.Range("J16").Value = 1.060950326364018
$ws.Range("K16").Value = 1.063657337972955
$ws.Range("L16").Value = 1.053804407909193
$ws.Range("M16").Value = 1.071109844569052
$ws.Range("N16").Value = 1.062456997520158

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.056447599266004
$ws.Range("D17").Value = 1.061904804836995
$ws.Range("E17").Value = 1.05191402735505
$ws.Range("F17").Value = 1.069471853225664
$ws.Range("I17").Value = 1.023843013862359
$ws.Range("J17").Value = 1.062494966149793
$ws.Range("K17").Value = 1.065182784420732
$ws.Range("L17").Value = 1.055225329233721
$ws.Range("M17").Value = 1.072725048573066
$ws.Range("N17").Value = 1.064003830871601

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.057427825519448
$ws.Range("D18").Value = 1.062834402903192
$ws.Range("E18").Value = 1.052783120706715
$ws.Range("F18").Value = 1.070453503568837
$ws.Range("I18").Value = 1.023826118601775
$ws.Range("J18").Value = 1.063391966919021
$ws.Range("K18").Value = 1.066068610634147
$ws.Range("L18").Value = 1.056050349092327
$ws.Range("M18").Value = 1.073663133312909
$ws.Range("N18").Value = 1.06490210548481

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.057761369119432
$ws.Range("D19").Value = 1.063150707248725
$ws.Range("E19").Value = 1.053078822368291
$ws.Range("F19").Value = 1.070787539717097
$ws.Range("I19").Value = 1.023820303822167
$ws.Range("J19").Value = 1.063697158033275
$ws.Range("K19").Value = 1.066369994910178
$ws.Range("L19").Value = 1.056331026795331
$ws.Range("M19").Value = 1.073982320603665
$ws.Range("N19").Value = 1.065207730005452

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.056266965082059
$ws.Range("D20").Value = 1.061733494672924
$ws.Range("E20").Value = 1.051753859935156
$ws.Range("F20").Value = 1.069290960420147
$ws.Range("I20").Value = 1.023846096030143
$ws.Range("J20").Value = 1.062329652819215
$ws.Range("K20").Value = 1.065019528200518
$ws.Range("L20").Value = 1.055073270745959
$ws.Range("M20").Value = 1.072552172197078
$ws.Range("N20").Value = 1.063838282777132

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.051359421885193
$ws.Range("D21").Value = 1.057078621464798
$ws.Range("E21").Value = 1.047400913506779
$ws.Range("F21").Value = 1.064376840916162
$ws.Range("I21").Value = 1.02392621804881
$ws.Range("J21").Value = 1.057836521481388
$ws.Range("K21").Value = 1.060582040038856
$ws.Range("L21").Value = 1.050939132355752
$ws.Range("M21").Value = 1.067854477015281
$ws.Range("N21").Value = 1.059338770677436

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048234518871111
$ws.Range("D22").Value = 1.054113985804673
$ws.Range("E22").Value = 1.044627753924189
$ws.Range("F22").Value = 1.061248198549688
$ws.Range("I22").Value = 1.023973741382756
$ws.Range("J22").Value = 1.054973718583547
$ws.Range("K22").Value = 1.057754430387521
$ws.Range("L22").Value = 1.048303833781941
$ws.Range("M22").Value = 1.064862307429474
$ws.Range("N22").Value = 1.056471902271112

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.049894930753249
$ws.Range("D23").Value = 1.055689297867532
$ws.Range("E23").Value = 1.046101398555662
$ws.Range("F23").Value = 1.062910554219863
$ws.Range("I23").Value = 1.023948818288664
$ws.Range("J23").Value = 1.056495031783047
$ws.Range("K23").Value = 1.059257066119123
$ws.Range("L23").Value = 1.049704365278724
$ws.Range("M23").Value = 1.0664522750748
$ws.Range("N23").Value = 1.057995375909853

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.056348598586832
$ws.Range("D24").Value = 1.061810914615462
$ws.Range("E24").Value = 1.051826244434646
$ws.Range("F24").Value = 1.069372710656522
$ws.Range("I24").Value = 1.023844704313072
$ws.Range("J24").Value = 1.062404363018574
$ws.Range("K24").Value = 1.065093308816833
$ws.Range("L24").Value = 1.055141991091542
$ws.Range("M24").Value = 1.072630300041993
$ws.Range("N24").Value = 1.063913099073544

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.063677354202725
$ws.Range("D25").Value = 1.068759876566058
$ws.Range("E25").Value = 1.058321289480979
$ws.Range("F25").Value = 1.0767129398829
$ws.Range("I25").Value = 1.02371141742819
$ws.Range("J25").Value = 1.069107373762433
$ws.Range("K25").Value = 1.071712301888705
$ws.Range("L25").Value = 1.061304683107327
$ws.Range("M25").Value = 1.079642208938693
$ws.Range("N25").Value = 1.070625628861505
